# "hoan thien sq diagram" - add a second slide (a "Section Break" style
# divider) right after the existing slide, carrying just the empty
# description/body placeholder (idx 12) from its layout.

$p = $ppt.ActivePresentation

# slideLayout26 ("Section Break") is the 12th CustomLayout that hangs off
# the presentation's 2nd slide master ("Office Theme" design) - it is the
# layout that owns a body placeholder with idx="12".
$design = $p.Designs.Item(2)
$layout = $design.SlideMaster.CustomLayouts.Item(12)

# Append the new slide (becomes slide 2) using that layout.
$slide = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

# The layout also seeds a title placeholder, but the final slide only
# keeps the body/description placeholder - drop the title shape.
$slide.Shapes.Item(1).Delete()

# Rename + clear the remaining placeholder's text (it stays empty).
$ph = $slide.Shapes.Item(1)
$ph.Name = "Chỗ dành sẵn cho Văn bản 2"
$ph.TextFrame.TextRange.Text = ""
